$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 875
$ws1.Range("F3").Value = 1450
$ws1.Range("F4").Value = 1107
$ws1.Range("F5").Value = 524
$ws1.Range("F6").Value = 225
$ws1.Range("F8").Value = 676
$ws1.Range("F9").Value = 252
$ws1.Range("F11").Value = 88
$ws1.Range("F14").Value = 2370
$ws1.Range("F15").Value = 434
$ws1.Range("F17").Value = 502
$ws1.Range("F20").Value = 114
$ws1.Range("F22").Value = 665
$ws1.Range("F25").Value = 964
$ws1.Range("F27").Value = 1574
$ws1.Range("F28").Value = 313
$ws1.Range("F29").Value = 28

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 122
$ws2.Range("F5").Value = 219

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 875
$ws4.Range("F4").Value = 1450
$ws4.Range("F5").Value = 1107
$ws4.Range("F6").Value = 122
$ws4.Range("F8").Value = 524
$ws4.Range("F9").Value = 225
$ws4.Range("F11").Value = 676
$ws4.Range("F13").Value = 252
$ws4.Range("F15").Value = 88
$ws4.Range("F18").Value = 2371
$ws4.Range("F19").Value = 219
$ws4.Range("F20").Value = 434
$ws4.Range("F22").Value = 502
$ws4.Range("F26").Value = 114
$ws4.Range("F31").Value = 665
$ws4.Range("F38").Value = 964
$ws4.Range("F40").Value = 1574
$ws4.Range("F41").Value = 313
$ws4.Range("F42").Value = 28
